# Update the marksheet's correct/total marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: correct-answer score value (B11): 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: total marks obtained (B12): 45 -> 75
$ws.Range("B12").Value = 75

# "Total" row: correct/total marks ratio string (E12): "45/84" -> "75/140"
$ws.Range("E12").Value = "75/140"
